$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Row 1 Col 1: 35÷7= -> 68÷3=
$cell = $tbl.Cell(1, 1)
$cell.Range.Text = "68÷3="
# Row 1 Col 2: 97÷8= -> 46÷6=
$cell = $tbl.Cell(1, 2)
$cell.Range.Text = "46÷6="
# Row 1 Col 3: 11÷6= -> 35÷3=
$cell = $tbl.Cell(1, 3)
$cell.Range.Text = "35÷3="
# Row 1 Col 4: 75÷9= -> 57÷6=
$cell = $tbl.Cell(1, 4)
$cell.Range.Text = "57÷6="
# Row 1 Col 5: 34÷5= -> 33÷6=
$cell = $tbl.Cell(1, 5)
$cell.Range.Text = "33÷6="
# Row 5 Col 1: 59÷4= -> 20÷8=
$cell = $tbl.Cell(5, 1)
$cell.Range.Text = "20÷8="
# Row 5 Col 2: 47÷6= -> 14÷9=
$cell = $tbl.Cell(5, 2)
$cell.Range.Text = "14÷9="
# Row 5 Col 3: 78÷3= -> 69÷8=
$cell = $tbl.Cell(5, 3)
$cell.Range.Text = "69÷8="
# Row 5 Col 4: 22÷2= -> 52÷9=
$cell = $tbl.Cell(5, 4)
$cell.Range.Text = "52÷9="
# Row 5 Col 5: 84÷8= -> 55÷5=
$cell = $tbl.Cell(5, 5)
$cell.Range.Text = "55÷5="
# Row 9 Col 1: 51÷2= -> 12÷5=
$cell = $tbl.Cell(9, 1)
$cell.Range.Text = "12÷5="
# Row 9 Col 2: 90÷8= -> 74÷2=
$cell = $tbl.Cell(9, 2)
$cell.Range.Text = "74÷2="
# Row 9 Col 3: 69÷2= -> 67÷8=
$cell = $tbl.Cell(9, 3)
$cell.Range.Text = "67÷8="
# Row 9 Col 4: 35÷3= -> 21÷8=
$cell = $tbl.Cell(9, 4)
$cell.Range.Text = "21÷8="
# Row 9 Col 5: 59÷2= -> 98÷3=
$cell = $tbl.Cell(9, 5)
$cell.Range.Text = "98÷3="
# Row 13 Col 1: 13÷7= -> 97÷8=
$cell = $tbl.Cell(13, 1)
$cell.Range.Text = "97÷8="
# Row 13 Col 2: 71÷8= -> 11÷6=
$cell = $tbl.Cell(13, 2)
$cell.Range.Text = "11÷6="
# Row 13 Col 3: 67÷3= -> 46÷7=
$cell = $tbl.Cell(13, 3)
$cell.Range.Text = "46÷7="
# Row 13 Col 4: 17÷9= -> 35÷4=
$cell = $tbl.Cell(13, 4)
$cell.Range.Text = "35÷4="
# Row 13 Col 5: 63÷2= -> 46÷3=
$cell = $tbl.Cell(13, 5)
$cell.Range.Text = "46÷3="
# Row 17 Col 1: 32÷2= -> 15÷5=
$cell = $tbl.Cell(17, 1)
$cell.Range.Text = "15÷5="
# Row 17 Col 2: 99÷8= -> 25÷8=
$cell = $tbl.Cell(17, 2)
$cell.Range.Text = "25÷8="
# Row 17 Col 3: 92÷6= -> 66÷3=
$cell = $tbl.Cell(17, 3)
$cell.Range.Text = "66÷3="
# Row 17 Col 4: 51÷6= -> 19÷7=
$cell = $tbl.Cell(17, 4)
$cell.Range.Text = "19÷7="
# Row 17 Col 5: 18÷7= -> 79÷4=
$cell = $tbl.Cell(17, 5)
$cell.Range.Text = "79÷4="
